$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.6862
$ws.Range("J3").Value = 3.0631
$ws.Range("J4").Value = 2.2513
$ws.Range("J5").Value = 2.9341
$ws.Range("J6").Value = 1.5682
$ws.Range("J7").Value = 2.2149
$ws.Range("J8").Value = 3.0488
$ws.Range("J9").Value = 2.391
$ws.Range("J10").Value = 3.3008
